# Auto-generated Excel COM-interop edit script
# Updates Leve market-price columns (H-N) for specific rows across the ALC, ARM,
# BSM, CRP, CUL, GSM, LTW and WVR sheets, matching a scheduled market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 111 (ALC)
$ws.Range("H111").Value = 10382.25
$ws.Range("I111").Value = 12843
$ws.Range("J111").Value = 3000
$ws.Range("K111").Value = 38529
$ws.Range("L111").Value = 9000
$ws.Range("M111").Value = -35462
$ws.Range("N111").Value = -15134

# Row 113 (ALC)
$ws.Range("H113").Value = 3389.9285
$ws.Range("I113").Value = 2425
$ws.Range("J113").Value = 3550.75
$ws.Range("K113").Value = 2425
$ws.Range("L113").Value = 3550.75
$ws.Range("M113").Value = 829
$ws.Range("N113").Value = -10058.75

$ws = $wb.Worksheets.Item("ARM")
# Row 2 (ARM)
$ws.Range("H2").Value = 1767.3158
$ws.Range("I2").Value = 1930.5
$ws.Range("J2").Value = 1487.5714
$ws.Range("K2").Value = 1930.5
$ws.Range("L2").Value = 1487.5714
$ws.Range("M2").Value = -1817.5
$ws.Range("N2").Value = -1713.5714

# Row 32 (ARM)
$ws.Range("H32").Value = 5350.104
$ws.Range("I32").Value = 4377.811
$ws.Range("J32").Value = 29333.334
$ws.Range("K32").Value = 4377.811
$ws.Range("L32").Value = 29333.334
$ws.Range("M32").Value = -4090.811
$ws.Range("N32").Value = -29907.334

# Row 61 (ARM)
$ws.Range("H61").Value = 5714
$ws.Range("I61").Value = 4717
$ws.Range("J61").Value = 10449.75
$ws.Range("K61").Value = 4717
$ws.Range("L61").Value = 10449.75
$ws.Range("M61").Value = -4505
$ws.Range("N61").Value = -10873.75

# Row 74 (ARM)
$ws.Range("H74").Value = 4863.5557
$ws.Range("I74").Value = 2728.8
$ws.Range("J74").Value = 9715.272000000001
$ws.Range("K74").Value = 2728.8
$ws.Range("L74").Value = 9715.272000000001
$ws.Range("M74").Value = -1854.8
$ws.Range("N74").Value = -11463.272

# Row 77 (ARM)
$ws.Range("H77").Value = 4863.5557
$ws.Range("I77").Value = 2728.8
$ws.Range("J77").Value = 9715.272000000001
$ws.Range("K77").Value = 13644
$ws.Range("L77").Value = 48576.36
$ws.Range("M77").Value = -9276
$ws.Range("N77").Value = -57312.36

# Row 116 (ARM)
$ws.Range("H116").Value = 1767.3158
$ws.Range("I116").Value = 1930.5
$ws.Range("J116").Value = 1487.5714
$ws.Range("K116").Value = 1930.5
$ws.Range("L116").Value = 1487.5714
$ws.Range("M116").Value = 363.5
$ws.Range("N116").Value = -6075.5714

# Row 122 (ARM)
$ws.Range("H122").Value = 2923.5625
$ws.Range("I122").Value = 5144.8
$ws.Range("J122").Value = 1913.909
$ws.Range("K122").Value = 15434.4
$ws.Range("L122").Value = 5741.727000000001
$ws.Range("M122").Value = -12984.4
$ws.Range("N122").Value = -10641.727

# Row 132 (ARM)
$ws.Range("H132").Value = 5441.027
$ws.Range("I132").Value = 1957.1852
$ws.Range("J132").Value = 14847.4
$ws.Range("K132").Value = 5871.5556
$ws.Range("L132").Value = 44542.2
$ws.Range("M132").Value = -3341.5556
$ws.Range("N132").Value = -49602.2

# Row 136 (ARM)
$ws.Range("H136").Value = 5714
$ws.Range("I136").Value = 4717
$ws.Range("J136").Value = 10449.75
$ws.Range("K136").Value = 14151
$ws.Range("L136").Value = 31349.25
$ws.Range("M136").Value = -11601
$ws.Range("N136").Value = -36449.25

$ws = $wb.Worksheets.Item("BSM")
# Row 3 (BSM)
$ws.Range("H3").Value = 1767.3158
$ws.Range("I3").Value = 1930.5
$ws.Range("J3").Value = 1487.5714
$ws.Range("K3").Value = 1930.5
$ws.Range("L3").Value = 1487.5714
$ws.Range("M3").Value = -1816.5
$ws.Range("N3").Value = -1715.5714

# Row 107 (BSM)
$ws.Range("H107").Value = 1868.619
$ws.Range("I107").Value = 1696.5294
$ws.Range("J107").Value = 2600
$ws.Range("K107").Value = 1696.5294
$ws.Range("L107").Value = 2600
$ws.Range("M107").Value = 223.4706000000001
$ws.Range("N107").Value = -6440

# Row 110 (BSM)
$ws.Range("H110").Value = 39702
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 39702
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 39702
$ws.Range("N110").Value = -47882

# Row 117 (BSM)
$ws.Range("H117").Value = 77300
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 77300
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 77300
$ws.Range("N117").Value = -86478

# Row 140 (BSM)
$ws.Range("H140").Value = 45974.785
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 45974.785
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 45974.785
$ws.Range("N140").Value = -56334.785

$ws = $wb.Worksheets.Item("CRP")
# Row 102 (CRP)
$ws.Range("H102").Value = 38000
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 38000
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 38000
$ws.Range("N102").Value = -42868

# Row 104 (CRP)
$ws.Range("H104").Value = 45000
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 45000
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 45000
$ws.Range("N104").Value = -50242

# Row 107 (CRP)
$ws.Range("H107").Value = 376.37036
$ws.Range("I107").Value = 366.36365
$ws.Range("J107").Value = 420.4
$ws.Range("K107").Value = 366.36365
$ws.Range("L107").Value = 420.4
$ws.Range("M107").Value = 1553.63635
$ws.Range("N107").Value = -4260.4

# Row 109 (CRP)
$ws.Range("H109").Value = 33333.332
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 33333.332
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 33333.332
$ws.Range("N109").Value = -35413.332

# Row 122 (CRP)
$ws.Range("H122").Value = 9384.214
$ws.Range("I122").Value = 3028.8823
$ws.Range("J122").Value = 19206.092
$ws.Range("K122").Value = 9086.6469
$ws.Range("L122").Value = 57618.276
$ws.Range("M122").Value = -6636.6469
$ws.Range("N122").Value = -62518.276

$ws = $wb.Worksheets.Item("CUL")
# Row 113 (CUL)
$ws.Range("H113").Value = 743.78
$ws.Range("I113").Value = 760.1829
$ws.Range("J113").Value = 669.05554
$ws.Range("K113").Value = 2280.5487
$ws.Range("L113").Value = 2007.16662
$ws.Range("M113").Value = -110.5487000000003
$ws.Range("N113").Value = -6347.16662

$ws = $wb.Worksheets.Item("GSM")
# Row 102 (GSM)
$ws.Range("H102").Value = 15012
$ws.Range("I102").Value = 15012
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 15012
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -13390
$ws.Range("N102").ClearContents()

# Row 107 (GSM)
$ws.Range("H107").Value = 345.7143
$ws.Range("I107").Value = 123.4
$ws.Range("J107").Value = 901.5
$ws.Range("K107").Value = 123.4
$ws.Range("L107").Value = 901.5
$ws.Range("M107").Value = 1796.6
$ws.Range("N107").Value = -4741.5

# Row 113 (GSM)
$ws.Range("H113").Value = 5099.4
$ws.Range("I113").Value = 6250
$ws.Range("J113").Value = 4332.3335
$ws.Range("K113").Value = 6250
$ws.Range("L113").Value = 4332.3335
$ws.Range("M113").Value = -4080
$ws.Range("N113").Value = -8672.333500000001

# Row 122 (GSM)
$ws.Range("I122").Value = 6150.5835
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 18451.7505
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -16001.7505
$ws.Range("N122").Value = -9400

$ws = $wb.Worksheets.Item("LTW")
# Row 61 (LTW)
$ws.Range("H61").Value = 14914.167
$ws.Range("I61").Value = 16360
$ws.Range("J61").Value = 9853.75
$ws.Range("K61").Value = 16360
$ws.Range("L61").Value = 9853.75
$ws.Range("M61").Value = -16158
$ws.Range("N61").Value = -10257.75

# Row 102 (LTW)
$ws.Range("H102").Value = 57000
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 57000
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 57000
$ws.Range("N102").Value = -63490

# Row 113 (LTW)
$ws.Range("H113").Value = 14914.167
$ws.Range("I113").Value = 16360
$ws.Range("J113").Value = 9853.75
$ws.Range("K113").Value = 16360
$ws.Range("L113").Value = 9853.75
$ws.Range("M113").Value = -14190
$ws.Range("N113").Value = -14193.75

# Row 122 (LTW)
$ws.Range("H122").Value = 4597.206
$ws.Range("I122").Value = 4562.5
$ws.Range("J122").Value = 4680.5
$ws.Range("K122").Value = 13687.5
$ws.Range("L122").Value = 14041.5
$ws.Range("M122").Value = -11237.5
$ws.Range("N122").Value = -18941.5

# Row 136 (LTW)
$ws.Range("H136").Value = 5414.6875
$ws.Range("I136").Value = 3974.4
$ws.Range("J136").Value = 6980.2173
$ws.Range("K136").Value = 11923.2
$ws.Range("L136").Value = 20940.6519
$ws.Range("M136").Value = -9373.200000000001
$ws.Range("N136").Value = -26040.6519

$ws = $wb.Worksheets.Item("WVR")
# Row 37 (WVR)
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

# Row 102 (WVR)
$ws.Range("H102").Value = 53000
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 53000
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 53000
$ws.Range("N102").Value = -59490

# Row 109 (WVR)
$ws.Range("H109").Value = 62900
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 62900
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 62900
$ws.Range("N109").Value = -65674

# Row 122 (WVR)
$ws.Range("H122").Value = 3755.5217
$ws.Range("I122").Value = 2461.375
$ws.Range("J122").Value = 6713.5713
$ws.Range("K122").Value = 7384.125
$ws.Range("L122").Value = 20140.7139
$ws.Range("M122").Value = -4934.125
$ws.Range("N122").Value = -25040.7139

# Row 136 (WVR)
$ws.Range("H136").Value = 3952.3333
$ws.Range("I136").Value = 2005.6945
$ws.Range("J136").Value = 6288.3
$ws.Range("K136").Value = 6017.083500000001
$ws.Range("L136").Value = 18864.9
$ws.Range("M136").Value = -3467.083500000001
$ws.Range("N136").Value = -23964.9
